# Daily attendance processing - 2025-11-01 12:35:03
# Reorders the "Recorded By" (column G) values so that "System"/"system"
# is moved from the front of the comma-separated list to the end,
# except where the other entry is admin@admin.com (left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -notlike "System,*") { continue }
    if ($val -like "*admin@admin.com*") { continue }

    $parts = $val -split ", "
    $newParts = $parts[1..($parts.Length - 1)] + $parts[0]
    $newVal = $newParts -join ", "

    $cell.Value = $newVal
}
